$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row block (row 53, which has
# the same B/C/D/E style pattern: s=5,5,1,1) onto the five new rows 57-61.
$ws.Range("B53:E53").Copy() | Out-Null
$ws.Range("B57:E61").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 57 ---------------------------------------------------------------
# Fill E, D, C (in that order) then B, matching the order new strings were
# first introduced into the shared-string table.
$ws.Range("E57").Value = "0:44`ncomeçando sessão de validação dos campos. Para isso é necessário trabalhar com Bean Validation ... requer dependencia:`n<!-- TOMAR CUIIDADO COM VERSÃO DO SPRING BOOT, A DEPENDENCIA ABAIXO`nATENDE AS VERSÕES DO SPRING ACIMA DE 2.3.0 -->`n<dependency>`n <groupId>org.springframework.boot</groupId>`n <artifactId>spring-boot-starter-validation</artifactId>`n</dependency>"
$ws.Range("D57").Value = "`n57. Validando formulários de cargo e departamento"
$ws.Range("C57").Value = "11. Validação Back-End"
$ws.Range("B57").Value = 57

# --- Row 58 ---------------------------------------------------------------
$ws.Range("E58").Value = "1:12`nanotações utilizadas no domain para validação dos seus atributos:`n@NotBlank(message=`"`")`n@Size()"
$ws.Range("D58").Value = "`n57. Validando formulários de cargo e departamento"
$ws.Range("C58").Value = "11. Validação Back-End"
$ws.Range("B58").Value = 57

# --- Row 60 (filled before row 59/61, per original authoring order) -------
$ws.Range("E60").Value = "6:33`nanotação @Valid nos métodos salvar e editar para validar campos entre controller e view HTML"
$ws.Range("D60").Value = "`n57. Validando formulários de cargo e departamento"
$ws.Range("C60").Value = "11. Validação Back-End"
$ws.Range("B60").Value = 57

# --- Row 61 -----------------------------------------------------------------
$ws.Range("E61").Value = "6:56`nadicionado um parametro do tipo BindingResult do spring para trabalhar em conjunto com a anotação @Valid nas validaçoes dos campos"
$ws.Range("D61").Value = "`n57. Validando formulários de cargo e departamento"
$ws.Range("C61").Value = "11. Validação Back-End"
$ws.Range("B61").Value = 57

# --- Row 59 -----------------------------------------------------------------
$ws.Range("E59").Value = "5:36`nIMPORTANTE: PARA FUNCIONAR DEVE ESTAR DENTRO DE UMA TAG `"FORM`" NO HTML caso contrário não funcionará.`ninclusão de fragmento de validação de campos nos forms de cadastro de departamento e de cargos para mostrar mensagens na tela.`n<div th:replace=`"fragments/validacao :: validacao`"></div>"
$ws.Range("D59").Value = "`n57. Validando formulários de cargo e departamento"
$ws.Range("C59").Value = "11. Validação Back-End"
$ws.Range("B59").Value = 57

# --- Row heights, matching the authored layout -----------------------------
$ws.Range("B57").EntireRow.RowHeight = 165
$ws.Range("B58").EntireRow.RowHeight = 75
$ws.Range("B59").EntireRow.RowHeight = 90
$ws.Range("B60").EntireRow.RowHeight = 45
$ws.Range("B61").EntireRow.RowHeight = 45

# --- Update the visible selection / scroll position to match the saved view -
$ws.Range("E60").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
